# Auto-generated PowerShell COM script to apply Goblin_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 84.71429
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H28").Value = 7262.5
$ws.Range("I28").Value = 7262.5
$ws.Range("K28").Value = 7262.5
$ws.Range("M28").Value = -6777.5
$ws.Range("H32").Value = 100000
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H33").Value = 1096.5385
$ws.Range("I33").Value = 241.44444
$ws.Range("K33").Value = 241.44444
$ws.Range("M33").Value = -12.44443999999999
$ws.Range("H39").Value = 211
$ws.Range("J39").Value = 1099.5
$ws.Range("L39").Value = 3298.5
$ws.Range("N39").Value = -3890.5
$ws.Range("H51").Value = 2230.7693
$ws.Range("I51").Value = 2187.5
$ws.Range("J51").Value = 2300
$ws.Range("K51").Value = 2187.5
$ws.Range("L51").Value = 2300
$ws.Range("M51").Value = -1703.5
$ws.Range("N51").Value = -3268
$ws.Range("H62").Value = 55125.438
$ws.Range("I62").Value = 74728
$ws.Range("K62").Value = 74728
$ws.Range("M62").Value = -74104
$ws.Range("H65").Value = 55125.438
$ws.Range("I65").Value = 74728
$ws.Range("K65").Value = 373640
$ws.Range("M65").Value = -370520
$ws.Range("H76").Value = 3092.8572
$ws.Range("I76").Value = 2997.5
$ws.Range("K76").Value = 2997.5
$ws.Range("M76").Value = -2682.5
$ws.Range("H79").Value = 3092.8572
$ws.Range("I79").Value = 2997.5
$ws.Range("K79").Value = 2997.5
$ws.Range("M79").Value = -1905.5
$ws.Range("H86").Value = 2483.182
$ws.Range("I86").Value = 2196.7144
$ws.Range("J86").Value = 2984.5
$ws.Range("K86").Value = 2196.7144
$ws.Range("L86").Value = 2984.5
$ws.Range("M86").Value = -1073.7144
$ws.Range("N86").Value = -5230.5
$ws.Range("H89").Value = 2483.182
$ws.Range("I89").Value = 2196.7144
$ws.Range("J89").Value = 2984.5
$ws.Range("K89").Value = 10983.572
$ws.Range("L89").Value = 14922.5
$ws.Range("M89").Value = -5367.572
$ws.Range("N89").Value = -26154.5
$ws.Range("H95").Value = 31232.875
$ws.Range("J95").Value = 31232.875
$ws.Range("L95").Value = 31232.875
$ws.Range("N95").Value = -36724.875
$ws.Range("H100").Value = 6663.3125
$ws.Range("I100").Value = 4499.8335
$ws.Range("K100").Value = 4499.8335
$ws.Range("M100").Value = -3958.8335
$ws.Range("H112").Value = 1833.2307
$ws.Range("J112").Value = 2524.3333
$ws.Range("L112").Value = 7572.999899999999
$ws.Range("N112").Value = -9788.999899999999
$ws.Range("H133").Value = 65000
$ws.Range("J133").Value = 65000
$ws.Range("L133").Value = 65000
$ws.Range("N133").Value = -75120
$ws.Range("H136").Value = 54246
$ws.Range("J136").Value = 54246
$ws.Range("L136").Value = 54246
$ws.Range("N136").Value = -64446
$ws.Range("H138").Value = 1663.0889
$ws.Range("J138").Value = 4484.9
$ws.Range("L138").Value = 13454.7
$ws.Range("N138").Value = -23734.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 20000
$ws.Range("J24").Value = 20000
$ws.Range("L24").Value = 20000
$ws.Range("N24").Value = -20748
$ws.Range("H45").Value = 726
$ws.Range("I45").Value = 726
$ws.Range("K45").Value = 726
$ws.Range("M45").Value = -349
$ws.Range("H97").Value = 1873.8
$ws.Range("I97").Value = 1830.4615
$ws.Range("J97").Value = 2155.5
$ws.Range("K97").Value = 1830.4615
$ws.Range("L97").Value = 2155.5
$ws.Range("M97").Value = -1334.4615
$ws.Range("N97").Value = -3147.5
$ws.Range("H100").Value = 20000
$ws.Range("J100").Value = 20000
$ws.Range("L100").Value = 20000
$ws.Range("N100").Value = -22164
$ws.Range("H109").Value = 46250
$ws.Range("J109").Value = 46250
$ws.Range("L109").Value = 46250
$ws.Range("N109").Value = -49024
$ws.Range("H132").Value = 2497.3845
$ws.Range("I132").Value = 2454.3333
$ws.Range("K132").Value = 7362.999899999999
$ws.Range("M132").Value = -4832.999899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 34495
$ws.Range("J76").Value = 34495
$ws.Range("L76").Value = 34495
$ws.Range("N76").Value = -35125
$ws.Range("H79").Value = 34495
$ws.Range("J79").Value = 34495
$ws.Range("L79").Value = 34495
$ws.Range("N79").Value = -36679
$ws.Range("H86").Value = 577581.9
$ws.Range("I86").Value = 2220.2083
$ws.Range("J86").Value = 2111879.5
$ws.Range("K86").Value = 2220.2083
$ws.Range("L86").Value = 2111879.5
$ws.Range("M86").Value = -1097.2083
$ws.Range("N86").Value = -2114125.5
$ws.Range("H89").Value = 577581.9
$ws.Range("I89").Value = 2220.2083
$ws.Range("J89").Value = 2111879.5
$ws.Range("K89").Value = 11101.0415
$ws.Range("L89").Value = 10559397.5
$ws.Range("M89").Value = -5485.041499999999
$ws.Range("N89").Value = -10570629.5
$ws.Range("H99").Value = 3590.7273
$ws.Range("I99").Value = 2624.75
$ws.Range("K99").Value = 2624.75
$ws.Range("M99").Value = -1126.75
$ws.Range("H105").Value = 2502.4849
$ws.Range("I105").Value = 2037.3448
$ws.Range("K105").Value = 2037.3448
$ws.Range("M105").Value = -290.3448000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 9799.8
$ws.Range("I25").Value = 4000
$ws.Range("J25").Value = 13666.333
$ws.Range("K25").Value = 4000
$ws.Range("L25").Value = 13666.333
$ws.Range("M25").Value = -3826
$ws.Range("N25").Value = -14014.333
$ws.Range("H99").Value = 3762.25
$ws.Range("I99").Value = 3762.25
$ws.Range("K99").Value = 3762.25
$ws.Range("M99").Value = -2264.25
$ws.Range("H107").Value = 999.75
$ws.Range("I107").Value = 999.75
$ws.Range("K107").Value = 999.75
$ws.Range("M107").Value = 920.25
$ws.Range("H126").Value = 3762.25
$ws.Range("I126").Value = 3762.25
$ws.Range("K126").Value = 11286.75
$ws.Range("M126").Value = -8816.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1000
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 3000
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -3346
$ws.Range("H140").Value = 1792
$ws.Range("I140").Value = 1516
$ws.Range("J140").Value = 4000
$ws.Range("K140").Value = 4548
$ws.Range("L140").Value = 12000
$ws.Range("M140").Value = 632
$ws.Range("N140").Value = -22360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 29999.666
$ws.Range("J49").Value = 29999.666
$ws.Range("L49").Value = 29999.666
$ws.Range("N49").Value = -30367.666
$ws.Range("H80").Value = 5272.0454
$ws.Range("I80").Value = 4320.5
$ws.Range("J80").Value = 6937.25
$ws.Range("K80").Value = 4320.5
$ws.Range("L80").Value = 6937.25
$ws.Range("M80").Value = -3322.5
$ws.Range("N80").Value = -8933.25
$ws.Range("H83").Value = 5272.0454
$ws.Range("I83").Value = 4320.5
$ws.Range("J83").Value = 6937.25
$ws.Range("K83").Value = 21602.5
$ws.Range("L83").Value = 34686.25
$ws.Range("M83").Value = -16610.5
$ws.Range("N83").Value = -44670.25
$ws.Range("H97").Value = 962
$ws.Range("I97").Value = 763.625
$ws.Range("K97").Value = 763.625
$ws.Range("M97").Value = -267.625
$ws.Range("H123").Value = 49499.25
$ws.Range("J123").Value = 49499.25
$ws.Range("L123").Value = 49499.25
$ws.Range("N123").Value = -54399.25
$ws.Range("H126").Value = 2693
$ws.Range("I126").Value = 2752.4614
$ws.Range("K126").Value = 8257.3842
$ws.Range("M126").Value = -5787.3842
$ws.Range("H132").Value = 2840.4443
$ws.Range("I132").Value = 2820
$ws.Range("J132").Value = 2856.8
$ws.Range("K132").Value = 8460
$ws.Range("L132").Value = 8570.400000000001
$ws.Range("M132").Value = -5930
$ws.Range("N132").Value = -13630.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3011.7646
$ws.Range("J7").Value = 6150
$ws.Range("L7").Value = 6150
$ws.Range("N7").Value = -6374
$ws.Range("H25").Value = 9966.667
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H46").Value = 3456.2144
$ws.Range("I46").Value = 2718.6
$ws.Range("K46").Value = 2718.6
$ws.Range("M46").Value = -2530.6
$ws.Range("H61").Value = 5214.6665
$ws.Range("I61").Value = 3038.5
$ws.Range("K61").Value = 3038.5
$ws.Range("M61").Value = -2836.5
$ws.Range("H82").Value = 2853.3076
$ws.Range("I82").Value = 1849.125
$ws.Range("J82").Value = 4460
$ws.Range("K82").Value = 1849.125
$ws.Range("L82").Value = 4460
$ws.Range("M82").Value = -1488.125
$ws.Range("N82").Value = -5182
$ws.Range("H85").Value = 2853.3076
$ws.Range("I85").Value = 1849.125
$ws.Range("J85").Value = 4460
$ws.Range("K85").Value = 1849.125
$ws.Range("L85").Value = 4460
$ws.Range("M85").Value = -601.125
$ws.Range("N85").Value = -6956
$ws.Range("H113").Value = 5214.6665
$ws.Range("I113").Value = 3038.5
$ws.Range("K113").Value = 3038.5
$ws.Range("M113").Value = -868.5
$ws.Range("H126").Value = 3011.7646
$ws.Range("J126").Value = 6150
$ws.Range("L126").Value = 18450
$ws.Range("N126").Value = -23390

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3284.3684
$ws.Range("J81").Value = 3494.8572
$ws.Range("L81").Value = 6989.7144
$ws.Range("N81").Value = -9111.7144
$ws.Range("H84").Value = 3284.3684
$ws.Range("J84").Value = 3494.8572
$ws.Range("L84").Value = 34948.572
$ws.Range("N84").Value = -45556.572

Write-Host "Applied all changes"